# Team19_Proposal.docx edit script
# Applies:
#  1. Cosmetic run-merge in "2. Al- Arafat ... 1420991042" line (no text change)
#  2. Cosmetic run-merge in "By using this project..." paragraph (no text change)
#  3. "We will use Hamcrest and JUnit framework..." -> add Mockito mention + new
#     sentence about Mockito being a mocking framework
#  4. "Framework: JUnit ... & Hamcrest" -> expand with Mockito/JUnit/Hamcrest version numbers
#  5. Move the _GoBack bookmark from after "5" to the blank paragraph after the
#     Framework/Hamcrest line

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "\t  2. Al- Arafat\t\t -  \t1420991042" - merge redundant runs (no text
#    change, just normalizes run boundaries the way Word would on re-save)
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(16)
$p1start = $p1.Range.Start
$delRange1 = $d.Range($p1start + 18, $p1start + 23)
$delRange1.Delete()
$insPoint1 = $d.Range($p1start + 18, $p1start + 18)
$insPoint1.InsertBefore(" -  " + [char]9)

# ---------------------------------------------------------------------------
# 2) "By using this project, admin will be able to ... saving it manually"
#    merge the 3 middle runs into one (no text change)
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(28)
$p2start = $p2.Range.Start
$delRange2 = $d.Range($p2start + 1, $p2start + 458)
$delRange2.Delete()
$insPoint2 = $d.Range($p2start + 1, $p2start + 1)
$insPoint2.InsertBefore('By using this project, admin will be able to create new bank account, add money or         deposit money, check balance of the accounts, check number of account holder, check percentage of interest and update it, search account for info, check for loan availability, give loan, withdrawing balance if account-holder wants and successful transaction of money. It will make a bank easier to perform these operations to save info rather than saving it manually')
$newRange2 = $d.Range($p2start + 1, $p2start + 458)
$newRange2.Font.Name = "Times New Roman"
$newRange2.Font.Size = 16

Write-Output "Step 1-2 done"

# ---------------------------------------------------------------------------
# 3) "We will use Hamcrest and JUnit framework in order to test whole
#    project. Hamcrest is a framework..." paragraph: mention Mockito too, and
#    add a new sentence describing Mockito as a mocking framework.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(42)
$p3start = $p3.Range.Start

# Replace the span covering "Hamcrest and JUnit framework in order to test
# whole project. " (old offsets 12-73 relative to paragraph start) with the
# expanded text.
$delRange3 = $d.Range($p3start + 12, $p3start + 73)
$delRange3.Delete()
$insPoint3 = $d.Range($p3start + 12, $p3start + 12)
$insPoint3.InsertBefore('Mockito, Hamcrest and JUnit framework in order to test whole project. Mockito is a mocking framework that tastes really good. It lets you write beautiful tests with a clean & simple API. ')

# Segment formatting (offsets relative to $p3start + 12):
#   0-7    "Mockito"                                   Times New Roman
#   7-9    ", "                                        Times New Roman
#   9-17   "Hamcrest"                                  Times New Roman
#   17-18  " "                                         Times New Roman
#   18-28  "and JUnit "                                Times New Roman
#   28-68  "framework in order to test whole project"  Times New Roman
#   68-69  "."                                         (cs)Times New Roman, black
#   69-70  " "                                         (cs)Arial, black
#   70-77  "Mockito"                                   (cs)Arial, black
#   77-185 " is a mocking framework ... clean & simple API" (cs)Arial, black
#   185-186 "."                                        (cs)Arial, black
#   186-187 " "                                        Times New Roman, black
$base3 = $p3start + 12

$seg = $d.Range($base3 + 68, $base3 + 69)
$seg.Font.NameBi = "Times New Roman"
$seg.Font.Color = 0

$seg = $d.Range($base3 + 69, $base3 + 185)
$seg.Font.NameBi = "Arial"
$seg.Font.Color = 0

$seg = $d.Range($base3 + 185, $base3 + 187)
$seg.Font.NameBi = "Arial"
$seg.Font.Color = 0

$seg = $d.Range($base3 + 186, $base3 + 187)
$seg.Font.NameAscii = "Times New Roman"
$seg.Font.NameBi = "Times New Roman"
$seg.Font.Color = 0

Write-Output "Step 3 done"
